$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 3 so that the existing masterDataPath row (old row 3)
# becomes row 6, and all following rows shift down by 3 (old row 4 -> new row 7, etc.)
$ws.Rows("3:5").Insert()

# Fill in the new "databasePath" row (new row 4) first so its strings are appended
# to the shared-strings table ahead of the "dbType" row's strings.
$ws.Range("A4").Value = "databasePath"
$ws.Range("B4").Value = "Path"
$ws.Range("C4").Value = "Dropbox/MDSTrungThien/Shared/Invenage.sqlite"

# Repurpose the old masterDataPath row (now row 3) into the new "dbType" row.
$ws.Range("A3").Value = "dbType"
$ws.Range("B3").Value = "Param"
$ws.Range("C3").Value = "SQLite"

# New "sqlUserName" / "sqlPassword" rows (rows 5 and 6). The row names are entered
# before the row values so the shared-strings table lists sqlUserName/sqlPassword
# ahead of umdocc/mdstrungthien.
$ws.Range("A5").Value = "sqlUserName"
$ws.Range("B5").Value = "Param"
$ws.Range("A6").Value = "sqlPassword"
$ws.Range("B6").Value = "Param"
$ws.Range("C5").Value = "umdocc"
$ws.Range("C6").Value = "mdstrungthien"

# Update the active cell selection to A4
$ws.Range("A4").Select()
